$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.341.34"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "1.667.12"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  +0.92%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.42"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5346"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2663"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06389"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.86"
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07850"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.556"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "1.681.67"
$ws.Range("E13").Value = "  +2.09%  "
$ws.Range("D14").Value = "1.895.39"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5545"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").Value = "0.0₅8196"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.95"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "26.361.75"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.673"
$ws.Range("E20").Value = "  +2.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.07"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.27"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.045"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.11"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1226"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.217"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.11"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05859"
$ws.Range("E30").Value = "  +0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.284"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.586"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9707"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.830"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5825"
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01609"
$ws.Range("D40").Value = "1.071.63"
$ws.Range("E40").Value = "  +4.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8643"
$ws.Range("E41").Value = "  +1.91%  "
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.22"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").Value = "1.805.35"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "0.0₈105"
$ws.Range("E48").Value = "  -6.00%  "
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.059"
$ws.Range("E50").Value = "  +2.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05164"
$ws.Range("E51").Value = "  +0.41%  "
